$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency data (prices and volume %) scraped on
# Mon Jun 10 01:12:03 UTC 2024. Some rows also swap rank order (B/C/D/E).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.545.52'
$ws.Range("E2").Value = '  +0.42%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.689.48'
$ws.Range("E3").Value = '  +0.20%  '

$ws.Range("E4").Value = '  +0.21%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '672.12'

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '160.89'
$ws.Range("E6").Value = '  +1.57%  '

$ws.Range("E7").Value = '  +0.09%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.498'
$ws.Range("E8").Value = '  +0.76%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.146'
$ws.Range("E9").Value = '  -0.62%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.14'
$ws.Range("E10").Value = '  +2.01%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.442'
$ws.Range("E11").Value = '  +0.97%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000234'
$ws.Range("E12").Value = '  +0.46%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '32.92'
$ws.Range("E13").Value = '  +2.03%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.685.64'
$ws.Range("E14").Value = '  +0.19%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '69.582.81'
$ws.Range("E15").Value = '  +0.44%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.116'
$ws.Range("E16").Value = '  +1.35%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '16.18'
$ws.Range("E17").Value = '  +1.03%  '

$ws.Range("E18").Value = '  +1.63%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '473.91'
$ws.Range("E19").Value = '  +0.62%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '9.77'
$ws.Range("E20").Value = '  -2.50%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.650'
$ws.Range("E21").Value = '  -0.14%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '80.37'
$ws.Range("E22").Value = '  +0.50%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.841.42'
$ws.Range("E23").Value = '  +0.38%  '

$ws.Range("B24").Value = 'Dai'
$ws.Range("C24").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.00'
$ws.Range("E24").Value = '  -0.03%  '

$ws.Range("B25").Value = 'PEPE'
$ws.Range("C25").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000126'
$ws.Range("E25").Value = '  +3.62%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.99'
$ws.Range("E26").Value = '  +0.64%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.17'
$ws.Range("E27").Value = '  +0.52%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.69'
$ws.Range("E28").Value = '  -0.41%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.72'
$ws.Range("E29").Value = '  -1.13%  '

$ws.Range("B30").Value = 'ImmutableX'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.02'
$ws.Range("E30").Value = '  +1.44%  '

$ws.Range("B31").Value = 'Kaspa'
$ws.Range("C31").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.168'
$ws.Range("E31").Value = '  +4.17%  '

$ws.Range("B32").Value = 'Binance-PegBSC-USD'
$ws.Range("C32").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.00'
$ws.Range("E32").Value = '  +0.06%  '

$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '26.87'
$ws.Range("E33").Value = '  +0.00%  '

$ws.Range("B34").Value = 'NEARProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.52'
$ws.Range("E34").Value = '  -2.23%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.689.81'
$ws.Range("E35").Value = '  +0.78%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '8.51'
$ws.Range("E36").Value = '  +3.98%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.12'
$ws.Range("E37").Value = '  -0.28%  '

$ws.Range("E38").Value = '  -0.05%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.25'
$ws.Range("E39").Value = '  -0.12%  '

$ws.Range("E40").Value = '  -0.07%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0908'
$ws.Range("E41").Value = '  +0.24%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '175.00'
$ws.Range("E42").Value = '  +2.13%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.937'
$ws.Range("E43").Value = '  -0.37%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '46.94'
$ws.Range("E44").Value = '  -1.19%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.75'
$ws.Range("E45").Value = '  +1.43%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.000276'
$ws.Range("E46").Value = '  -2.08%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.29'
$ws.Range("E47").Value = '  +0.25%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '27.65'
$ws.Range("E48").Value = '  +0.17%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.08'
$ws.Range("E49").Value = '  -1.33%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.87'
$ws.Range("E50").Value = '  +1.07%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.266'
$ws.Range("E51").Value = '  -0.54%  '
